$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '64.605.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -0.19%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.144.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -0.28%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'" + '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.12%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '576.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +0.72%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '148.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -1.58%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.03%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '3.142.71'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -0.23%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '0.526'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.40%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -2.45%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '6.13'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -0.90%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '0.500'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -0.68%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '0.0000262'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +2.37%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '37.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -2.23%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '3.657.27'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -0.36%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '64.655.57'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -0.33%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = "'" + '7.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -1.33%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = "'" + '3.144.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.48%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +0.37%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '503.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -2.84%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '14.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -0.78%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = "'" + '15.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.11%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = "'" + '0.712'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -3.59%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '7.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -1.87%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '84.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -1.25%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -0.07%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '2.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -1.26%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '8.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +1.23%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -1.04%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '2.80'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +5.25%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '27.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -1.55%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -0.01%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +1.02%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '6.18'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +0.85%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '6.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -2.10%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '54.51'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -2.19%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '0.0888'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +2.71%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '476.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -2.34%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -2.16%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '2.93'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -2.47%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '8.70'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +0.52%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '2.994.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -3.85%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -4.58%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = "'" + '2.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -2.80%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = "'" + '0.281'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -5.87%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '28.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -3.91%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.0₃0580'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +0.34%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -0.03%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -1.70%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -2.95%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = "'" + '2.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +11.79%  '
$ws.Range('E51').Style = 'Normal'
